$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: name -> "vanilla", price -> 200
$ws.Range("B2").Value = "vanilla"
$ws.Range("C2").Value = 200

# Remove row 3 entirely (it contained Id=2, Name=shit, Price=300)
$ws.Rows(3).Delete()

# Set width of column C to match the new layout
$ws.Columns(3).ColumnWidth = 14

# Update the selected cell to C3 (now the first empty row below the data)
$ws.Range("C3").Select()
